$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124; existing rows 124:216 shift down to 125:217
$ws.Rows("124:124").Insert()

# Populate the newly inserted row 124 with the new weekly data point
$ws.Cells.Item(124, 1).Value = 11
$ws.Cells.Item(124, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(124, 3).Value = "Bíobío"
$ws.Cells.Item(124, 4).Value = 45090
$ws.Cells.Item(124, 5).Value = 8
$ws.Cells.Item(124, 6).Value = 100112043
$ws.Cells.Item(124, 7).Value = "Pepino ensalada"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 100
$ws.Cells.Item(124, 11).Value = 11000
$ws.Cells.Item(124, 12).Value = 12000
$ws.Cells.Item(124, 13).Value = 11500
$ws.Cells.Item(124, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(124, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(124, 16).Value = 192
$ws.Cells.Item(124, 17).Value = 60
$ws.Cells.Item(124, 18).Value = "Hortaliza"
